# Bugfixed evaluation and simulated rt_data for components.
# The naive AR2 forecast vectors shifted by one period (a new, earlier
# observation was inserted) and the recomputed forecast ratios (columns
# C and E) changed for (almost) every existing row. A brand-new final
# row (53) was also appended with the latest observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 53 is new; copy formatting (date style, s="2") from row 52 column A
# so the new date cell A53 keeps the same number format as the rest of
# column A before we overwrite its value.
$ws.Range("A52").Copy()
$ws.Range("A53").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 2
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 1.75539628881467
$ws.Range("D2").Value = 2008

# Row 3
$ws.Range("A3").Value = 39583

# Row 4
$ws.Range("A4").Value = 39765
$ws.Range("B4").Value = 2008
$ws.Range("C4").Value = 2.213911448916162
$ws.Range("D4").Value = 2009

# Row 5
$ws.Range("A5").Value = 39948
$ws.Range("C5").ClearContents()

# Row 6
$ws.Range("A6").Value = 40130
$ws.Range("B6").Value = 2009
$ws.Range("C6").Value = 2.533533936850563
$ws.Range("D6").Value = 2010

# Row 7
$ws.Range("A7").Value = 40310
$ws.Range("C7").Value = 2.208165160720954
$ws.Range("E7").Value = 1.903751357432193

# Row 8
$ws.Range("A8").Value = 40494
$ws.Range("B8").Value = 2010
$ws.Range("C8").Value = 2.088987486264915
$ws.Range("D8").Value = 2011
$ws.Range("E8").Value = 1.485473821631844

# Row 9
$ws.Range("A9").Value = 40676
$ws.Range("C9").Value = 1.614140618728332
$ws.Range("E9").Value = 1.770808585446004

# Row 10
$ws.Range("A10").Value = 40862
$ws.Range("B10").Value = 2011
$ws.Range("C10").Value = 1.212544822741002
$ws.Range("D10").Value = 2012
$ws.Range("E10").Value = 1.799394172339341

# Row 11
$ws.Range("A11").Value = 41044
$ws.Range("C11").Value = 1.625793900975747
$ws.Range("E11").Value = 1.586821460965226

# Row 12
$ws.Range("A12").Value = 41228
$ws.Range("B12").Value = 2012
$ws.Range("C12").Value = 1.196776590518644
$ws.Range("D12").Value = 2013
$ws.Range("E12").Value = 1.2151583353186

# Row 13
$ws.Range("A13").Value = 41409
$ws.Range("C13").Value = 0.8049364973309325
$ws.Range("E13").Value = 1.421244400332

# Row 14
$ws.Range("A14").Value = 41592
$ws.Range("B14").Value = 2013
$ws.Range("C14").Value = 0.4712609263772594
$ws.Range("D14").Value = 2014
$ws.Range("E14").Value = 1.107727073902187

# Row 15
$ws.Range("A15").Value = 41774
$ws.Range("C15").Value = 0.5775251578155283
$ws.Range("E15").Value = 1.341244385861273

# Row 16
$ws.Range("A16").Value = 41957
$ws.Range("B16").Value = 2014
$ws.Range("C16").Value = 0.8783377572271434
$ws.Range("D16").Value = 2015
$ws.Range("E16").Value = 1.612081704302182

# Row 17
$ws.Range("A17").Value = 42137
$ws.Range("C17").Value = 1.901826580533572
$ws.Range("E17").Value = 1.53605963063923

# Row 18
$ws.Range("A18").Value = 42321
$ws.Range("B18").Value = 2015
$ws.Range("C18").Value = 2.29066283401107
$ws.Range("D18").Value = 2016
$ws.Range("E18").Value = 2.221748592150097

# Row 19
$ws.Range("A19").Value = 42503
$ws.Range("C19").Value = 2.590339257583607
$ws.Range("E19").Value = 1.672072534917302

# Row 20
$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = 2016
$ws.Range("C20").Value = 4.109890522944348
$ws.Range("D20").Value = 2017
$ws.Range("E20").Value = 2.932944072183674

# Row 21
$ws.Range("A21").Value = 42867
$ws.Range("C21").Value = 1.713587272940131
$ws.Range("E21").Value = 1.721854626734953

# Row 22
$ws.Range("A22").Value = 43053
$ws.Range("B22").Value = 2017
$ws.Range("C22").Value = 1.336316831462692
$ws.Range("D22").Value = 2018
$ws.Range("E22").Value = 1.104283769064729

# Row 23
$ws.Range("A23").Value = 43145
$ws.Range("C23").Value = 1.808022822788802
$ws.Range("E23").Value = 1.867774135387434

# Row 24
$ws.Range("A24").Value = 43235
$ws.Range("C24").Value = 1.05432456490544
$ws.Range("E24").Value = 1.415552619392124

# Row 25
$ws.Range("A25").Value = 43326
$ws.Range("C25").Value = 1.299469465444592
$ws.Range("E25").Value = 1.645976944955962

# Row 26
$ws.Range("A26").Value = 43418
$ws.Range("B26").Value = 2018
$ws.Range("C26").Value = 1.197912858979611
$ws.Range("D26").Value = 2019
$ws.Range("E26").Value = 1.649865498505276

# Row 27
$ws.Range("A27").Value = 43510
$ws.Range("C27").Value = 2.247656020455691
$ws.Range("E27").Value = 1.934175841213626

# Row 28
$ws.Range("A28").Value = 43600
$ws.Range("C28").Value = 1.566023898188384
$ws.Range("E28").Value = 1.644188696416427

# Row 29
$ws.Range("A29").Value = 43691
$ws.Range("C29").Value = 1.75655962297816
$ws.Range("E29").Value = 2.050351917667315

# Row 30
$ws.Range("A30").Value = 43783
$ws.Range("B30").Value = 2019
$ws.Range("C30").Value = 1.727537197898665
$ws.Range("D30").Value = 2020
$ws.Range("E30").Value = 2.284828905445169

# Row 31
$ws.Range("A31").Value = 43875
$ws.Range("C31").Value = 2.684967757027334
$ws.Range("E31").Value = 2.075491449101596

# Row 32
$ws.Range("A32").Value = 43966
$ws.Range("C32").Value = 2.155932165770968
$ws.Range("E32").Value = 1.805141163113122

# Row 33
$ws.Range("A33").Value = 44068
$ws.Range("C33").Value = 2.980209378995857
$ws.Range("E33").Value = 2.653391228709334

# Row 34
$ws.Range("A34").Value = 44159
$ws.Range("B34").Value = 2020
$ws.Range("C34").Value = 3.647228437274408
$ws.Range("D34").Value = 2021
$ws.Range("E34").Value = 3.474365686630398

# Row 35
$ws.Range("A35").Value = 44251
$ws.Range("C35").Value = 1.906805170974435
$ws.Range("E35").Value = 1.938263709207333

# Row 36
$ws.Range("A36").Value = 44341
$ws.Range("C36").Value = 2.443967114785739
$ws.Range("E36").Value = 2.026008136667135

# Row 37
$ws.Range("A37").Value = 44432
$ws.Range("C37").Value = 1.954146674711188
$ws.Range("E37").Value = 1.667670056759474

# Row 38
$ws.Range("A38").Value = 44525
$ws.Range("B38").Value = 2021
$ws.Range("C38").Value = 2.777797690741424
$ws.Range("D38").Value = 2022
$ws.Range("E38").Value = 1.742844348069261

# Row 39
$ws.Range("A39").Value = 44617
$ws.Range("C39").Value = 1.467237762893392
$ws.Range("E39").Value = 1.862063279188941

# Row 40
$ws.Range("A40").Value = 44706
$ws.Range("C40").Value = 0.388123216496683
$ws.Range("E40").Value = 1.819907598678561

# Row 41
$ws.Range("A41").Value = 44798
$ws.Range("C41").Value = 2.69102598245059
$ws.Range("E41").Value = 3.239034933968399

# Row 42
$ws.Range("A42").Value = 44890
$ws.Range("B42").Value = 2022
$ws.Range("C42").Value = 0.6994919452575576
$ws.Range("D42").Value = 2023
$ws.Range("E42").Value = 0.5651273241891186

# Row 43
$ws.Range("A43").Value = 44981
$ws.Range("C43").Value = 0.2809429127725194
$ws.Range("E43").Value = 1.721404396148163

# Row 44
$ws.Range("A44").Value = 45071
$ws.Range("C44").Value = -2.811030211656218
$ws.Range("E44").Value = 0.8407670860975047

# Row 45
$ws.Range("A45").Value = 45163
$ws.Range("C45").Value = -1.669605379075589
$ws.Range("E45").Value = 0.6305126186323617

# Row 46
$ws.Range("A46").Value = 45254
$ws.Range("B46").Value = 2023
$ws.Range("C46").Value = -1.432689847121871
$ws.Range("D46").Value = 2024
$ws.Range("E46").Value = 0.4518870186319468

# Row 47
$ws.Range("A47").Value = 45345
$ws.Range("C47").Value = 1.069839250900739
$ws.Range("E47").Value = 1.634674340565567

# Row 48
$ws.Range("A48").Value = 45436
$ws.Range("C48").Value = 1.250641979737566
$ws.Range("E48").Value = 1.466559393695466

# Row 49
$ws.Range("A49").Value = 45534
$ws.Range("C49").Value = 1.780300968358017
$ws.Range("E49").Value = 1.762346671645298

# Row 50
$ws.Range("A50").Value = 45618
$ws.Range("B50").Value = 2024
$ws.Range("C50").Value = 2.033479419175133
$ws.Range("D50").Value = 2025
$ws.Range("E50").Value = 1.959987726090251

# Row 51
$ws.Range("A51").Value = 45713
$ws.Range("C51").Value = 2.97447584856072
$ws.Range("E51").Value = 1.953339169714385

# Row 52
$ws.Range("A52").Value = 45800
$ws.Range("C52").Value = 2.302179720973463
$ws.Range("E52").Value = 1.805984941845473

# Row 53
$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 2025
$ws.Range("C53").Value = 2.481068287768839
$ws.Range("D53").Value = 2026
$ws.Range("E53").Value = 1.908500198348873
